$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 36; all existing rows 36.. shift down by one.
$ws.Rows.Item(36).Insert()

# The data that used to live in row 36 is now in row 37. Duplicate it back into
# the newly-inserted row 36, then overwrite the Fecha (D) and Volumen (M) cells
# with the new values for this observation.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(36, $col).Value = $ws.Cells.Item(37, $col).Value2
}

$ws.Cells.Item(36, 4).Value = 44575   # D36 Fecha
$ws.Cells.Item(36, 13).Value = 250    # M36 Volumen
